# Update column G (K = strikeouts) with regenerated values.
# Previously these held "Strike#" based values; now regenerated to the K stat.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = 4
    3 = 9
    4 = 6
    5 = 11
    6 = 6
    7 = 3
    8 = 9
    9 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
